$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1931809676459194
$ws.Range("C3").Value = 0.3996566649694033
$ws.Range("C4").Value = 0.4602448435348362
$ws.Range("C5").Value = 0.5321724539119475
$ws.Range("C6").Value = 0.7579028661353813
$ws.Range("C7").Value = 1.253562338231302
$ws.Range("C8").Value = 0.8660741495955956
$ws.Range("C9").Value = 0.1818605364817396
$ws.Range("C10").Value = 0.07977418509674408
$ws.Range("C11").Value = -0.2970238499474278
$ws.Range("C12").Value = -0.507806226476132
$ws.Range("C13").Value = -0.5515046536471334
$ws.Range("C14").Value = 0.02404072806993659
$ws.Range("C15").Value = 0.2057539375803512
$ws.Range("C16").Value = -0.0555252746796395
$ws.Range("C17").Value = -0.1842593512688773
$ws.Range("C18").Value = -0.117190169769105
$ws.Range("C19").Value = -0.9219526200415645
$ws.Range("C20").Value = -1.027631261042042
$ws.Range("C21").Value = -0.994799733359641
$ws.Range("C22").Value = -1.357523728040164
$ws.Range("C23").Value = -1.436363496252893
$ws.Range("C24").Value = -1.384375511964711
$ws.Range("C25").Value = -1.208523948690917
$ws.Range("C26").Value = -2.203514422569788
$ws.Range("C27").Value = -2.25383584010028
$ws.Range("C28").Value = -1.694002501593443
$ws.Range("C29").Value = -0.7388195486982237
$ws.Range("C30").Value = -1.1353957149983
$ws.Range("C31").Value = -1.15191447423667
$ws.Range("C32").Value = -0.8765168772969435
$ws.Range("C33").Value = -0.7341883454077257
$ws.Range("C34").Value = -0.5542359184481118
$ws.Range("C35").Value = -0.7170935443645891
$ws.Range("C36").Value = -0.7163759421450732
$ws.Range("C37").Value = -0.8078134098658957
$ws.Range("C38").Value = -0.2074126431067658
$ws.Range("C39").Value = -0.1816222815636901
$ws.Range("C40").Value = -0.58485141349193
$ws.Range("C41").Value = -1.136309853296368
$ws.Range("C42").Value = -0.4716908278625516
$ws.Range("C43").Value = 0.5869115600794291
$ws.Range("C44").Value = 0.742580528021594
$ws.Range("C45").Value = 1.090717274284225
$ws.Range("C46").Value = 1.141012838903114
$ws.Range("C47").Value = 1.107826745279047
$ws.Range("C48").Value = 1.180343725096971
$ws.Range("C49").Value = 1.1413292403118
$ws.Range("C50").Value = 0.7454337975086591
$ws.Range("C51").Value = 0.5852354106138273
$ws.Range("C52").Value = 0.7466603069525758
$ws.Range("C53").Value = 1.385864701512679
$ws.Range("C54").Value = 1.333886277167564
$ws.Range("C55").Value = 0.9582665668632715
$ws.Range("C56").Value = 0.8566876594023194
$ws.Range("C57").Value = 0.3059084382741258
$ws.Range("C58").Value = 0.3293463138101435
$ws.Range("C59").Value = 1.016700246231457
$ws.Range("C60").Value = 1.063000856164264
$ws.Range("C61").Value = 1.019861959598145
$ws.Range("C62").Value = 1.044107195861401
$ws.Range("C63").Value = 1.955712778352741
$ws.Range("C64").Value = 1.938492711667397
$ws.Range("C65").Value = 1.853472242177145
$ws.Range("C66").Value = 1.55571227004196
$ws.Range("C67").Value = 1.032854664406683
$ws.Range("C68").Value = 1.128884167742598
$ws.Range("C69").Value = 1.231123923857274
$ws.Range("C70").Value = 1.577427547516803
$ws.Range("C71").Value = 1.333298080051702
$ws.Range("C72").Value = 1.255306145889133
$ws.Range("C73").Value = 1.311906613095918
$ws.Range("C74").Value = 1.290595733548666
$ws.Range("C75").Value = 0.4004417789580117
$ws.Range("C76").Value = 0.2098346370413737
$ws.Range("C77").Value = -0.134488317849601
$ws.Range("C78").Value = 0.5115738473705508
$ws.Range("C79").Value = 0.7680553565753787
$ws.Range("C80").Value = 0.6736363256059597
$ws.Range("C81").Value = 0.2671018265170643
$ws.Range("C82").Value = 0.3824097813653705
$ws.Range("C83").Value = -0.05124223606150702
$ws.Range("C84").Value = -0.1677363021246468
$ws.Range("C85").Value = -0.2330062845752818
$ws.Range("C86").Value = -0.3177632585026564
$ws.Range("C87").Value = -0.617171183702468
$ws.Range("C88").Value = -0.3239133254851859
$ws.Range("C89").Value = -0.1422659638130288
$ws.Range("C90").Value = -0.6549982355236162
$ws.Range("C91").Value = -0.7707967073906291
$ws.Range("C92").Value = -0.6935361792956022
$ws.Range("C93").Value = -0.2906846082282996
$ws.Range("C94").Value = -0.860610993712711
$ws.Range("C95").Value = -0.8081540746184284
$ws.Range("C96").Value = -0.5049407239096796
$ws.Range("C97").Value = -0.4757981464961759
$ws.Range("C98").Value = -0.4962077765769174
$ws.Range("C99").Value = -0.2629895660432023
$ws.Range("C100").Value = -0.4304122511906426
$ws.Range("C101").Value = -1.300689294146847
$ws.Range("C102").Value = -1.458790370815161
$ws.Range("C103").Value = -1.790776833134185
$ws.Range("C104").Value = -1.876661858279846
$ws.Range("C105").Value = -1.84153869254446
$ws.Range("C106").Value = -1.141406841527059
$ws.Range("C107").Value = -0.3176650160195481
$ws.Range("C108").Value = -0.555913106383176
$ws.Range("C109").Value = -0.42686865935449
$ws.Range("C110").Value = 0.4004517366521265
$ws.Range("C111").Value = 0.5921419271041274
$ws.Range("C112").Value = 0.6544450714398065
$ws.Range("C113").Value = 1.139324114415675
$ws.Range("C114").Value = 1.449840401084676
$ws.Range("C115").Value = 1.936202759884607
$ws.Range("C116").Value = 2.252205088517316
$ws.Range("C117").Value = 2.147487919592872
$ws.Range("C118").Value = 1.931769210205264
$ws.Range("C119").Value = 2.381948740921486
$ws.Range("C120").Value = 2.297596118668992
$ws.Range("C121").Value = 2.373618498258861
$ws.Range("C122").Value = 1.802260768733745
$ws.Range("C123").Value = 1.570652381522931
$ws.Range("C124").Value = 1.635283208114581
$ws.Range("C125").Value = 1.191080968177286
$ws.Range("C126").Value = 1.587705311991289
$ws.Range("C127").Value = 1.275578621991441
$ws.Range("C128").Value = 1.415091349802655
$ws.Range("C129").Value = 1.260165844998416
$ws.Range("C130").Value = 1.645182397911672
$ws.Range("C131").Value = 0.6839333083194281
$ws.Range("C132").Value = 0.623822310425387
$ws.Range("C133").Value = 0.6656792479090403
$ws.Range("C134").Value = 0.9124062851953783
$ws.Range("C135").Value = 0.4173659595475535
$ws.Range("C136").Value = 0.5767352406515757
$ws.Range("C137").Value = 1.086004056313912
